$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (C value, E value)
$updates = @{
    41  = @(113, 11056852)
    63  = @(14362, 36193311)
    64  = @(5213, 20442161)
    65  = @(2017, 13657300)
    70  = @(15734, 24684987)
    83  = @(3415, 115800208)
    91  = @(151149, 482480640)
    92  = @(409212, 1596566709)
    93  = @(209624, 1309639040)
    95  = @(50792, 933718057)
    96  = @(17307, 795758616)
    143 = @(64958, 373531948)
    146 = @(4269, 161502590)
    172 = @(22702, 44686673)
    184 = @(68737, 134191957)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}
